$wb = $excel.ActiveWorkbook

# Rename sheets to unify the DataNode naming convention
$wsProperty1 = $wb.Worksheets.Item("Property1")
$wsProperty1.Name = "DataNode_1"

$wsProperty2 = $wb.Worksheets.Item("Property2")
$wsProperty2.Name = "DataNode_2"

# Switch the active/selected tab to the second sheet (DataNode_2)
$wsProperty2.Activate()
